# Adding 'single' option to colour roller
# Adds a new "test" worksheet (after the last existing sheet) that holds the
# single-colour pick lists used by the roller: pens / pencils / other.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "test"

# Header row
$newSheet.Range("A1").Value = "pens"
$newSheet.Range("B1").Value = "pencils"
$newSheet.Range("C1").Value = "other"

# Column A (pens) filled top-to-bottom first
$newSheet.Range("A2").Value = "red"
$newSheet.Range("A3").Value = "orange"
$newSheet.Range("A4").Value = "yellow"
$newSheet.Range("A5").Value = "green"
$newSheet.Range("A6").Value = "blue"
$newSheet.Range("A7").Value = "indigo"
$newSheet.Range("A8").Value = "violet"

# Column B (pencils) filled next
$newSheet.Range("B2").Value = "black"
$newSheet.Range("B3").Value = "white"

# Column C (other) filled last
$newSheet.Range("C2").Value = "blender"

$newSheet.Range("O8").Select()
